$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Symmetrize the connectivity matrix: for each off-diagonal pair (i,j),
# the cell holding the smaller value is updated to match the larger value,
# matching the values supplied by the commit diff.

$ws.Range("B1").Value2 = 0.86931731463265716
$ws.Range("C1").Value2 = 0.93425937174813067
$ws.Range("AR2").Value2 = 0.97914190164217874
$ws.Range("AW2").Value2 = 0.84886673395194756
$ws.Range("B3").Value2 = 0.80288264433762824
$ws.Range("E3").Value2 = 0.60684333266411916
$ws.Range("B4").Value2 = 0.92225875438591054
$ws.Range("C4").Value2 = 0.92120186931090176
$ws.Range("W4").Value2 = 0.99880816297970043
$ws.Range("BB4").Value2 = 0.95491093808994687
$ws.Range("U5").Value2 = 0.5513298669920671
$ws.Range("G6").Value2 = 0.95796234799830082
$ws.Range("I7").Value2 = 0.90115629271469988
$ws.Range("BJ7").Value2 = 0.98548510018685909
$ws.Range("F8").Value2 = 0.82077355338093461
$ws.Range("H9").Value2 = 0.91779655196081444
$ws.Range("BL9").Value2 = 0.75265097011463433
$ws.Range("H10").Value2 = 0.58652174386306655
$ws.Range("I10").Value2 = 0.85930779592357776
$ws.Range("BM10").Value2 = 0.73331483409337261
$ws.Range("M11").Value2 = 0.8807913469314449
$ws.Range("K12").Value2 = 0.9429700404853858
$ws.Range("L13").Value2 = 0.60894465926553332
$ws.Range("O13").Value2 = 0.77452728822519523
$ws.Range("AT13").Value2 = 0.80816457664296204
$ws.Range("L14").Value2 = 0.86374338108871473
$ws.Range("P14").Value2 = 0.70449206785102825
$ws.Range("N15").Value2 = 0.98210630189892212
$ws.Range("P15").Value2 = 0.84889854254370534
$ws.Range("Q16").Value2 = 0.77349210399747637
$ws.Range("R16").Value2 = 0.81654111512456695
$ws.Range("O17").Value2 = 0.82166954488914823
$ws.Range("R17").Value2 = 0.64223696213247328
$ws.Range("AE17").Value2 = 0.97649941221101
$ws.Range("S18").Value2 = 0.70033176699586219
$ws.Range("Q19").Value2 = 0.93611581853218706
$ws.Range("R20").Value2 = 0.69471989103879106
$ws.Range("S20").Value2 = 0.98481964327290594
$ws.Range("U20").Value2 = 0.80499624871829334
$ws.Range("V20").Value2 = 0.61161717409620797
$ws.Range("W21").Value2 = 0.64134651306910473
$ws.Range("BP21").Value2 = 0.98186664466475138
$ws.Range("Y22").Value2 = 0.84782604440852827
$ws.Range("V23").Value2 = 0.72461886363713002
$ws.Range("X23").Value2 = 0.92152643453651006
$ws.Range("BD23").Value2 = 0.87481657035345406
$ws.Range("Z24").Value2 = 0.94465338214089156
$ws.Range("AE24").Value2 = 0.63905147326567224
$ws.Range("AG24").Value2 = 0.84408196771093513
$ws.Range("AA25").Value2 = 0.98963596850243563
$ws.Range("BO26").Value2 = 0.70885439023599806
$ws.Range("Z27").Value2 = 0.86079265148515938
$ws.Range("AB27").Value2 = 0.8178556284680909
$ws.Range("M28").Value2 = 0.79346755661962942
$ws.Range("U28").Value2 = 0.83042646486078486
$ws.Range("Z28").Value2 = 0.91649806085585217
$ws.Range("AD28").Value2 = 0.91691595679561488
$ws.Range("BG28").Value2 = 0.99315224001567226
$ws.Range("AA29").Value2 = 0.84049294064322633
$ws.Range("AD29").Value2 = 0.7972128563037848
$ws.Range("AE29").Value2 = 0.90108052169649211
$ws.Range("AF30").Value2 = 0.97314033180216231
$ws.Range("AD31").Value2 = 0.80142824902175147
$ws.Range("AN31").Value2 = 0.74247362542893636
$ws.Range("AE33").Value2 = 0.79455663736173643
$ws.Range("AF33").Value2 = 0.66435285760963214
$ws.Range("AF34").Value2 = 0.96155272705863581
$ws.Range("AI34").Value2 = 0.82682621006228429
$ws.Range("AJ34").Value2 = 0.87403592681138087
$ws.Range("AG35").Value2 = 0.89040769825166044
$ws.Range("AI36").Value2 = 0.70598842774481896
$ws.Range("AL36").Value2 = 0.99853768234198625
$ws.Range("AI37").Value2 = 0.64172019865496843
$ws.Range("AJ37").Value2 = 0.79151671989102645
$ws.Range("AK38").Value2 = 0.75654736195806749
$ws.Range("AN38").Value2 = 0.9033014066625431
$ws.Range("BC38").Value2 = 0.64649276319436522
$ws.Range("AK39").Value2 = 0.85850421307868463
$ws.Range("AO39").Value2 = 0.97086750846422121
$ws.Range("AM40").Value2 = 0.95296236911004395
$ws.Range("AO40").Value2 = 0.63107231766906935
$ws.Range("AA41").Value2 = 0.63547964853644778
$ws.Range("AP41").Value2 = 0.95102212420748322
$ws.Range("Y42").Value2 = 0.66596407646425981
$ws.Range("AQ42").Value2 = 0.86669862371510376
$ws.Range("AR42").Value2 = 0.93109576243455505
$ws.Range("AR43").Value2 = 0.6781833260925616
$ws.Range("AS44").Value2 = 0.93092830359561729
$ws.Range("S45").Value2 = 0.83741133776194343
$ws.Range("AS46").Value2 = 0.84542485641023668
$ws.Range("AS47").Value2 = 0.96523219314694386
$ws.Range("AT47").Value2 = 0.97901299612510306
$ws.Range("BF47").Value2 = 0.72108294577756771
$ws.Range("AT48").Value2 = 0.72242426031536211
$ws.Range("AU48").Value2 = 0.91444466491293563
$ws.Range("AX48").Value2 = 0.88362049110187613
$ws.Range("BA48").Value2 = 0.96565259189458463
$ws.Range("BE48").Value2 = 0.86844431324835003
$ws.Range("X49").Value2 = 0.86189220020772517
$ws.Range("AV49").Value2 = 0.87941771668975044
$ws.Range("J51").Value2 = 0.74612923805035758
$ws.Range("AJ51").Value2 = 0.70986339046161406
$ws.Range("AX51").Value2 = 0.75493340389955721
$ws.Range("AS52").Value2 = 0.86228735863425354
$ws.Range("AX52").Value2 = 0.93942789211027888
$ws.Range("AY52").Value2 = 0.89736767643153414
$ws.Range("AU53").Value2 = 0.99511788444721039
$ws.Range("AZ54").Value2 = 0.71186152836875549
$ws.Range("BB55").Value2 = 0.95287321740836317
$ws.Range("BD55").Value2 = 0.96769333295502058
$ws.Range("BO55").Value2 = 0.61571656348811576
$ws.Range("BB56").Value2 = 0.71024438263182876
$ws.Range("BE56").Value2 = 0.83364884431275554
$ws.Range("BC57").Value2 = 0.92805947633667885
$ws.Range("BD58").Value2 = 0.81333574323893343
$ws.Range("BE58").Value2 = 0.80129850640792566
$ws.Range("BJ58").Value2 = 0.94851904636538076
$ws.Range("AJ59").Value2 = 0.78873223343584997
$ws.Range("BF59").Value2 = 0.88672641982371048
$ws.Range("BJ60").Value2 = 0.73020632573576127
$ws.Range("BH61").Value2 = 0.79531078483046769
$ws.Range("N62").Value2 = 0.98198837237778958
$ws.Range("AT62").Value2 = 0.94512761948839019
$ws.Range("BI62").Value2 = 0.74073947110196303
$ws.Range("BI63").Value2 = 0.76454872628554904
$ws.Range("BL63").Value2 = 0.99921192093917299
$ws.Range("BM63").Value2 = 0.5442727475455722
$ws.Range("AQ65").Value2 = 0.82429493175786406
$ws.Range("AS65").Value2 = 0.94735908466948482
$ws.Range("BL65").Value2 = 0.97797979022975201
$ws.Range("BN65").Value2 = 0.94432026082823173
$ws.Range("J66").Value2 = 0.88054119256005969
$ws.Range("BL66").Value2 = 0.78069521176080414
$ws.Range("A68").Value2 = 0.79934061656435373
$ws.Range("E68").Value2 = 0.92022034291092281
$ws.Range("BN68").Value2 = 0.83010604949956623
